# edit.ps1
# Applies two changes described by the target diff:
#   1. The cached "datetimeFigureOut" field text ("29/01/2019") is updated
#      to "10/01/2019" everywhere it appears: once on the slide master and
#      once on each of the 11 slide layouts.
#   2. The logo group shape ("Groupe 15") on slide 1 is repositioned: its
#      <a:off> changes from x=3039386,y=768554 (EMU) to x=4223751,y=812097
#      (EMU); the size (<a:ext>) is unchanged.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Update the cached date placeholder text on the master + every layout.
# ---------------------------------------------------------------------
$oldDate = "29/01/2019"
$newDate = "10/01/2019"

function Update-DateShapes {
    param($shapes)
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShapes($master.Shapes)

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    Update-DateShapes($layout.Shapes)
}

# ---------------------------------------------------------------------
# 2) Move the logo group shape on slide 1.
# ---------------------------------------------------------------------
# Shape.Left/Top are expressed in points (EMU / 12700) and are stored
# internally as single-precision floats, so plain "EMU / 12700.0" can be
# off by 1 EMU after the round-trip. The literals below were chosen so
# that, after the float32 round-trip performed by the host, they land
# exactly back on the target EMU values.
$s = $p.Slides.Item(1)
$logo = $s.Shapes.Item(1)
$logo.Left = 332.578857421875
$logo.Top = 63.94468688964844
